$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "16÷2="
$t.Cell(1, 2).Range.Text = "90÷5="
$t.Cell(1, 3).Range.Text = "95÷8="
$t.Cell(1, 4).Range.Text = "69÷5="
$t.Cell(1, 5).Range.Text = "87÷8="
$t.Cell(5, 1).Range.Text = "53÷5="
$t.Cell(5, 2).Range.Text = "77÷7="
$t.Cell(5, 3).Range.Text = "37÷9="
$t.Cell(5, 4).Range.Text = "27÷8="
$t.Cell(5, 5).Range.Text = "21÷8="
$t.Cell(9, 1).Range.Text = "36÷3="
$t.Cell(9, 2).Range.Text = "32÷3="
$t.Cell(9, 3).Range.Text = "41÷4="
$t.Cell(9, 4).Range.Text = "75÷3="
$t.Cell(9, 5).Range.Text = "79÷8="
$t.Cell(13, 1).Range.Text = "59÷9="
$t.Cell(13, 2).Range.Text = "81÷7="
$t.Cell(13, 3).Range.Text = "20÷6="
$t.Cell(13, 4).Range.Text = "90÷9="
$t.Cell(13, 5).Range.Text = "70÷9="
$t.Cell(17, 1).Range.Text = "66÷8="
$t.Cell(17, 2).Range.Text = "63÷3="
$t.Cell(17, 3).Range.Text = "34÷8="
$t.Cell(17, 4).Range.Text = "72÷3="
$t.Cell(17, 5).Range.Text = "42÷4="
